$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAY_01")

$ws.Range("B3").Value = "27.Oct.2016"
$ws.Range("C3").Value = "Connect Spring Hibernate Application with Oracle"
$ws.Range("E3").Value = "http://o7planning.org/en/10305/simple-login-java-web-application-using-spring-mvc-spring-security-and-spring-jdbc                                                                                                                https://community.oracle.com/thread/2278282"
$ws.Range("F3").Value = "1 Hours"
$ws.Range("G3").Value = "N"

$ws.Hyperlinks.Add($ws.Range("E3"), "https://community.oracle.com/thread/2278282", "", "", "https://community.oracle.com/thread/2278282") | Out-Null

$ws.Range("G3").Select()
